# Update the "Notes" sheet (first sheet) of the uganda-igf-resources workbook.
#
# The sheet stores free-form documentation text, one line per row in column A.
# This edit:
#   - rewords the "Description:" line
#   - rewords the "Source:" line and adds a new "Source-link:" line after it
#   - rewords the license line and adds a new line pointing to the license info
#
# NOTE: on this particular (single-column) worksheet, assigning directly to an
# existing cell's .Value via Range("A<n>").Value = "..." has been observed to
# behave like an "insert row" (shifting everything below down one row and
# dropping the last row) rather than overwriting in place. To avoid that, we
# clear the sheet's existing rows first and then rewrite the full, final set
# of lines from row 1 downward in order, which is stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# Remove all existing rows of content on the Notes sheet.
$ws.Range("A1:A17").EntireRow.Delete()

$ws.Range("A1").Value = "Name: uganda-igf-resources"
$ws.Range("A2").Value = "Description: Locally Raised Revenues"
$ws.Range("A3").Value = "Units of measure: unit"
$ws.Range("A4").Value = "Source: Uganda budget information - Ministry of Finance, Planning and Economic Development"
$ws.Range("A5").Value = "Source-link: http://www.budget.go.ug/"
$ws.Range("A6").Value = ""
$ws.Range("A7").Value = "Notes:"
$ws.Range("A8").Value = ""
$ws.Range("A9").Value = "This data contains information that has been converted from current native currency units (NCU) to constant US Dollars. The NCU values are contained in the 'value-ncu' column, while the converted and deflated values are contained in the 'value' column."
$ws.Range("A10").Value = ""
$ws.Range("A11").Value = "On the 'Data-wide-value' sheet, we have provided the indicator in a wide format. The values you see listed there are from the 'value' column."
$ws.Range("A12").Value = ""
$ws.Range("A13").Value = ""
$ws.Range("A14").Value = ""
$ws.Range("A15").Value = "The following is data downloaded from Development Initiative's Datahub: http://devinit.org/data"
$ws.Range("A16").Value = "It is licensed under a Creative Commons Attribution 4.0 International license."
$ws.Range("A17").Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"
$ws.Range("A18").Value = "For concerns, questions, or corrections: please email info@devinit.org"
$ws.Range("A19").Value = "Copyright Development Initiatives Poverty Research Ltd. 2015"
